$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $escaped = $value -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

$excel.ScreenUpdating = $false

Set-TextValue $ws.Range('D2') '62.439.36'
Set-TextValue $ws.Range('E2') '  -0.91%  '

Set-TextValue $ws.Range('D3') '2.434.08'
Set-TextValue $ws.Range('E3') '  -0.52%  '

Set-TextValue $ws.Range('D5') '572.95'
Set-TextValue $ws.Range('E5') '  +0.37%  '

Set-TextValue $ws.Range('D6') '143.55'
Set-TextValue $ws.Range('E6') '  -2.22%  '

Set-TextValue $ws.Range('D8') '0.530'
Set-TextValue $ws.Range('E8') '  -0.53%  '

Set-TextValue $ws.Range('D9') '2.430.95'
Set-TextValue $ws.Range('E9') '  -0.74%  '

Set-TextValue $ws.Range('E10') '  -4.06%  '

Set-TextValue $ws.Range('E11') '  +0.74%  '

Set-TextValue $ws.Range('E12') '  -0.51%  '

Set-TextValue $ws.Range('E13') '  -1.49%  '

Set-TextValue $ws.Range('D14') '26.51'
Set-TextValue $ws.Range('E14') '  -1.33%  '

Set-TextValue $ws.Range('E15') '  -3.43%  '

Set-TextValue $ws.Range('D16') '2.871.74'
Set-TextValue $ws.Range('E16') '  -0.91%  '

Set-TextValue $ws.Range('D17') '62.295.03'
Set-TextValue $ws.Range('E17') '  -0.94%  '

Set-TextValue $ws.Range('D18') '2.440.11'
Set-TextValue $ws.Range('E18') '  -0.47%  '

Set-TextValue $ws.Range('E19') '  -2.68%  '

Set-TextValue $ws.Range('E20') '  -1.61%  '

Set-TextValue $ws.Range('D21') '325.77'
Set-TextValue $ws.Range('E21') '  +0.40%  '

Set-TextValue $ws.Range('E22') '  -0.80%  '

Set-TextValue $ws.Range('E23') '  +4.32%  '

Set-TextValue $ws.Range('E24') '  +0.17%  '

Set-TextValue $ws.Range('D25') '65.34'
Set-TextValue $ws.Range('E25') '  -1.38%  '

Set-TextValue $ws.Range('D26') '627.19'
Set-TextValue $ws.Range('E26') '  +1.19%  '

Set-TextValue $ws.Range('D27') '9.03'
Set-TextValue $ws.Range('E27') '  +4.60%  '

Set-TextValue $ws.Range('B28') 'PEPE'
Set-TextValue $ws.Range('C28') 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range('D28') '0.0₃0959'
Set-TextValue $ws.Range('E28') '  -6.60%  '

Set-TextValue $ws.Range('B29') 'WrappedeETH'
Set-TextValue $ws.Range('C29') 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue $ws.Range('D29') '2.560.59'
Set-TextValue $ws.Range('E29') '  -0.45%  '

Set-TextValue $ws.Range('B30') 'Binance-PegBSC-USD'
Set-TextValue $ws.Range('C30') 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range('D30') '1.01'
Set-TextValue $ws.Range('E30') '  +0.64%  '

Set-TextValue $ws.Range('B31') 'Fetch.AI'
Set-TextValue $ws.Range('C31') 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D31') '1.43'
Set-TextValue $ws.Range('E31') '  -3.49%  '

Set-TextValue $ws.Range('D32') '8.03'
Set-TextValue $ws.Range('E32') '  -2.76%  '

Set-TextValue $ws.Range('D33') '1.90'
Set-TextValue $ws.Range('E33') '  +0.37%  '

Set-TextValue $ws.Range('D34') '0.134'
Set-TextValue $ws.Range('E34') '  -6.40%  '

Set-TextValue $ws.Range('E35') '  -2.37%  '

Set-TextValue $ws.Range('E36') '  +0.30%  '

Set-TextValue $ws.Range('E37') '  -3.72%  '

Set-TextValue $ws.Range('E38') '  -1.97%  '

Set-TextValue $ws.Range('E39') '  -1.15%  '

Set-TextValue $ws.Range('B40') 'Monero'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D40') '146.89'
Set-TextValue $ws.Range('E40') '  +0.82%  '

Set-TextValue $ws.Range('B41') 'RenderToken'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Range('D41') '5.22'
Set-TextValue $ws.Range('E41') '  -3.30%  '

Set-TextValue $ws.Range('E42') '  -3.55%  '

Set-TextValue $ws.Range('D43') '42.28'
Set-TextValue $ws.Range('E43') '  +1.42%  '

Set-TextValue $ws.Range('E45') '  -5.51%  '

Set-TextValue $ws.Range('D46') '145.33'
Set-TextValue $ws.Range('E46') '  -1.32%  '

Set-TextValue $ws.Range('E47') '  -0.53%  '

Set-TextValue $ws.Range('D48') '0.0523'
Set-TextValue $ws.Range('E48') '  -3.17%  '

Set-TextValue $ws.Range('E49') '  -0.76%  '

Set-TextValue $ws.Range('D50') '19.72'
Set-TextValue $ws.Range('E50') '  -4.37%  '

Set-TextValue $ws.Range('E51') '  -1.89%  '

$excel.CutCopyMode = $false
$excel.ScreenUpdating = $true